$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "bool" example cell from "let flag = TRUE;" to "let flag = true;"
$ws.Range("C5").Value = "let flag = true;"

# Update the selection to C5 (was E5)
$ws.Range("C5").Select()
